$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 4587.778
$ws.Range("I38").Value = 882
$ws.Range("J38").Value = 11999.333
$ws.Range("K38").Value = 2646
$ws.Range("L38").Value = 35997.999
$ws.Range("M38").Value = -2274
$ws.Range("N38").Value = -36741.999
$ws.Range("H53").Value = 295.33334
$ws.Range("J53").Value = 229.71428
$ws.Range("L53").Value = 229.71428
$ws.Range("N53").Value = -1503.71428
$ws.Range("H80").Value = 1688.9445
$ws.Range("I80").Value = 3950.1667
$ws.Range("J80").Value = 558.3333
$ws.Range("K80").Value = 11850.5001
$ws.Range("L80").Value = 1674.9999
$ws.Range("M80").Value = -10852.5001
$ws.Range("N80").Value = -3670.9999
$ws.Range("H83").Value = 1688.9445
$ws.Range("I83").Value = 3950.1667
$ws.Range("J83").Value = 558.3333
$ws.Range("K83").Value = 35551.5003
$ws.Range("L83").Value = 5024.9997
$ws.Range("M83").Value = -30559.5003
$ws.Range("N83").Value = -15008.9997
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("N129").ClearContents()
$ws.Range("H137").Value = 3770.5
$ws.Range("I137").Value = 940.6923
$ws.Range("K137").Value = 2822.0769
$ws.Range("M137").Value = -272.0769
$ws.Range("H138").Value = 3264.2126
$ws.Range("J138").Value = 3259.6924
$ws.Range("L138").Value = 9779.0772
$ws.Range("N138").Value = -20059.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7365354
$ws.Range("I32").Value = 7818969.5
$ws.Range("K32").Value = 7818969.5
$ws.Range("M32").Value = -7818682.5
$ws.Range("H61").Value = 17903312
$ws.Range("I61").Value = 41672596
$ws.Range("J61").Value = 76347.625
$ws.Range("K61").Value = 41672596
$ws.Range("L61").Value = 76347.625
$ws.Range("M61").Value = -41672384
$ws.Range("N61").Value = -76771.625
$ws.Range("H74").Value = 7358583
$ws.Range("I74").Value = 8929701
$ws.Range("J74").Value = 26696.666
$ws.Range("K74").Value = 8929701
$ws.Range("L74").Value = 26696.666
$ws.Range("M74").Value = -8928827
$ws.Range("N74").Value = -28444.666
$ws.Range("H77").Value = 7358583
$ws.Range("I77").Value = 8929701
$ws.Range("J77").Value = 26696.666
$ws.Range("K77").Value = 44648505
$ws.Range("L77").Value = 133483.33
$ws.Range("M77").Value = -44644137
$ws.Range("N77").Value = -142219.33
$ws.Range("H136").Value = 17903312
$ws.Range("I136").Value = 41672596
$ws.Range("J136").Value = 76347.625
$ws.Range("K136").Value = 125017788
$ws.Range("L136").Value = 229042.875
$ws.Range("M136").Value = -125015238
$ws.Range("N136").Value = -234142.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 40000
$ws.Range("I32").Value = 40000
$ws.Range("K32").Value = 40000
$ws.Range("M32").Value = -39616
$ws.Range("H76").Value = 48399.8
$ws.Range("J76").Value = 48399.8
$ws.Range("L76").Value = 48399.8
$ws.Range("N76").Value = -49029.8
$ws.Range("H79").Value = 48399.8
$ws.Range("J79").Value = 48399.8
$ws.Range("L79").Value = 48399.8
$ws.Range("N79").Value = -50583.8
$ws.Range("H134").Value = 47889.26
$ws.Range("I134").Value = 4154.684
$ws.Range("K134").Value = 12464.052
$ws.Range("M134").Value = -9929.052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 842703.4399999999
$ws.Range("I31").Value = 26280.834
$ws.Range("K31").Value = 26280.834
$ws.Range("M31").Value = -25985.834
$ws.Range("H34").Value = 842703.4399999999
$ws.Range("I34").Value = 26280.834
$ws.Range("K34").Value = 26280.834
$ws.Range("M34").Value = -26078.834
$ws.Range("H107").Value = 1132.3
$ws.Range("I107").Value = 771.2857
$ws.Range("K107").Value = 771.2857
$ws.Range("M107").Value = 1148.7143

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 499.5
$ws.Range("J86").Value = 499
$ws.Range("L86").Value = 1497
$ws.Range("N86").Value = -3869
$ws.Range("H89").Value = 499.5
$ws.Range("J89").Value = 499
$ws.Range("L89").Value = 4491
$ws.Range("N89").Value = -16347
$ws.Range("H108").Value = 695.5
$ws.Range("I108").Value = 695.5
$ws.Range("K108").Value = 2086.5
$ws.Range("M108").Value = 793.5
$ws.Range("H141").Value = 112589.61
$ws.Range("I141").Value = 128854.664
$ws.Range("K141").Value = 386563.992
$ws.Range("M141").Value = -381383.992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 16401.7
$ws.Range("I43").Value = 9288.143
$ws.Range("J43").Value = 33000
$ws.Range("K43").Value = 9288.143
$ws.Range("L43").Value = 33000
$ws.Range("M43").Value = -9137.143
$ws.Range("N43").Value = -33302
$ws.Range("H57").Value = 22500
$ws.Range("J57").Value = 30000
$ws.Range("L57").Value = 30000
$ws.Range("N57").Value = -31640
$ws.Range("H113").Value = 3213.65
$ws.Range("I113").Value = 2376.6428
$ws.Range("J113").Value = 5166.6665
$ws.Range("K113").Value = 2376.6428
$ws.Range("L113").Value = 5166.6665
$ws.Range("M113").Value = -206.6428000000001
$ws.Range("N113").Value = -9506.666499999999
$ws.Range("H122").Value = 2999.8
$ws.Range("I122").Value = 2999.75
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 8999.25
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -6549.25
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1380.0435
$ws.Range("I16").Value = 1368.3572
$ws.Range("J16").Value = 1398.2222
$ws.Range("K16").Value = 1368.3572
$ws.Range("L16").Value = 1398.2222
$ws.Range("M16").Value = -1198.3572
$ws.Range("N16").Value = -1738.2222
$ws.Range("H100").Value = 9109.9
$ws.Range("I100").Value = 10633.333
$ws.Range("J100").Value = 8457
$ws.Range("K100").Value = 10633.333
$ws.Range("L100").Value = 8457
$ws.Range("M100").Value = -10092.333
$ws.Range("N100").Value = -9539
$ws.Range("H136").Value = 128400.53
$ws.Range("I136").Value = 86483.586
$ws.Range("J136").Value = 229001.2
$ws.Range("K136").Value = 259450.758
$ws.Range("L136").Value = 687003.6000000001
$ws.Range("M136").Value = -256900.758
$ws.Range("N136").Value = -692103.6000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 52000
$ws.Range("I42").Value = 52000
$ws.Range("K42").Value = 52000
$ws.Range("M42").Value = -51622
$ws.Range("H81").Value = 38894.09
$ws.Range("I81").Value = 24759.445
$ws.Range("J81").Value = 102500
$ws.Range("K81").Value = 49518.89
$ws.Range("L81").Value = 205000
$ws.Range("M81").Value = -48457.89
$ws.Range("N81").Value = -207122
$ws.Range("H84").Value = 38894.09
$ws.Range("I84").Value = 24759.445
$ws.Range("J84").Value = 102500
$ws.Range("K84").Value = 247594.45
$ws.Range("L84").Value = 1025000
$ws.Range("M84").Value = -242290.45
$ws.Range("N84").Value = -1035608
$ws.Range("H113").Value = 621.075
$ws.Range("I113").Value = 588.2963
$ws.Range("J113").Value = 689.1539
$ws.Range("K113").Value = 1764.8889
$ws.Range("L113").Value = 2067.4617
$ws.Range("M113").Value = 405.1111000000001
$ws.Range("N113").Value = -6407.4617
$ws.Range("H114").Value = 124995
$ws.Range("J114").Value = 124995
$ws.Range("L114").Value = 124995
$ws.Range("N114").Value = -133673
